$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Web_SIGNIN moves up from row 4 to row 3
$ws.Range("A3").Value = "Web_SIGNIN"
$ws.Range("B3").Value = 54
$ws.Range("C3").Value = "N"

# Row 4: Web_SEARCH moves up from row 5 to row 4 (no B value)
$ws.Range("A4").Value = "Web_SEARCH"
$ws.Range("B4").Clear()
$ws.Range("C4").Value = "N"

# Row 5: Web_MESSAGING moves up from row 6 to row 5 (no B value)
$ws.Range("A5").Value = "Web_MESSAGING"
$ws.Range("C5").Value = "N"

# Row 6: Web_QUICK_MESSAGES moves down from row 3 to row 6 (no B value)
$ws.Range("A6").Value = "Web_QUICK_MESSAGES"
$ws.Range("C6").Value = "Y"

# Update selection to C6
$ws.Range("C6").Select()

$wb.Save()
